$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows from 40 data rows (A2:E41) to 53 data rows (A2:E54).
# Append the extra rows needed at the bottom before rewriting all values.
$ws.Range("A42:A54").EntireRow.Insert()

$data = New-Object 'object[,]' 53,5
$data[0,0] = 1165327
$data[0,1] = 'loperamide Oral Liquid Product'
$data[0,2] = 'SCDG'
$data[0,3] = 'Active'
$data[0,4] = '[''loperamide'']'
$data[1,0] = 1165328
$data[1,1] = 'loperamide Oral Product'
$data[1,2] = 'SCDG'
$data[1,3] = 'Active'
$data[1,4] = '[''loperamide'']'
$data[2,0] = 1167781
$data[2,1] = 'Imodium Oral Liquid Product'
$data[2,2] = 'SBDG'
$data[2,3] = 'Active'
$data[2,4] = '[''loperamide'']'
$data[3,0] = 1167782
$data[3,1] = 'Imodium Oral Product'
$data[3,2] = 'SBDG'
$data[3,3] = 'Active'
$data[3,4] = '[''loperamide'']'
$data[4,0] = 1174579
$data[4,1] = 'Dolorex Solution Injectable Product'
$data[4,2] = 'SBDG'
$data[4,3] = 'Active'
$data[4,4] = '[''butorphanol'']'
$data[5,0] = 1176518
$data[5,1] = 'Diamode Oral Product'
$data[5,2] = 'SBDG'
$data[5,3] = 'Active'
$data[5,4] = '[''loperamide'']'
$data[6,0] = 1250685
$data[6,1] = 'loperamide hydrochloride 0.133 MG/ML Oral Suspension'
$data[6,2] = 'SCD'
$data[6,3] = 'Active'
$data[6,4] = '[''loperamide'']'
$data[7,0] = 1250693
$data[7,1] = 'loperamide hydrochloride 0.133 MG/ML Oral Suspension [Imodium]'
$data[7,2] = 'SBD'
$data[7,3] = 'Active'
$data[7,4] = '[''loperamide'']'
$data[8,0] = 1302739
$data[8,1] = 'butorphanol 10 MG/ML Injectable Solution'
$data[8,2] = 'SCD'
$data[8,3] = 'Active'
$data[8,4] = '[''butorphanol'']'
$data[9,0] = 1302741
$data[9,1] = 'butorphanol 10 MG/ML Injectable Solution [Dolorex Solution]'
$data[9,2] = 'SBD'
$data[9,3] = 'Active'
$data[9,4] = '[''butorphanol'']'
$data[10,0] = 1310925
$data[10,1] = 'butorphanol Injectable Solution [Butorphic]'
$data[10,2] = 'SBDF'
$data[10,3] = 'Active'
$data[10,4] = '[''butorphanol'']'
$data[11,0] = 1310926
$data[11,1] = 'Butorphic Injectable Product'
$data[11,2] = 'SBDG'
$data[11,3] = 'Active'
$data[11,4] = '[''butorphanol'']'
$data[12,0] = 1310927
$data[12,1] = 'butorphanol 10 MG/ML Injectable Solution [Butorphic]'
$data[12,2] = 'SBD'
$data[12,3] = 'Active'
$data[12,4] = '[''butorphanol'']'
$data[13,0] = 1489989
$data[13,1] = 'butorphanol Injectable Solution [Torbugesic]'
$data[13,2] = 'SBDF'
$data[13,3] = 'Active'
$data[13,4] = '[''butorphanol'']'
$data[14,0] = 1489990
$data[14,1] = 'Torbugesic Injectable Product'
$data[14,2] = 'SBDG'
$data[14,3] = 'Active'
$data[14,4] = '[''butorphanol'']'
$data[15,0] = 1489991
$data[15,1] = 'butorphanol 10 MG/ML Injectable Solution [Torbugesic]'
$data[15,2] = 'SBD'
$data[15,3] = 'Active'
$data[15,4] = '[''butorphanol'']'
$data[16,0] = 1594650
$data[16,1] = 'buprenorphine 1.8 MG/ML Injectable Solution'
$data[16,2] = 'SCD'
$data[16,3] = 'Active'
$data[16,4] = '[''buprenorphine'']'
$data[17,0] = 1594654
$data[17,1] = 'Simbadol Injectable Product'
$data[17,2] = 'SBDG'
$data[17,3] = 'Active'
$data[17,4] = '[''buprenorphine'']'
$data[18,0] = 1594655
$data[18,1] = 'buprenorphine 1.8 MG/ML Injectable Solution [Simbadol]'
$data[18,2] = 'SBD'
$data[18,3] = 'Active'
$data[18,4] = '[''buprenorphine'']'
$data[19,0] = 1809204
$data[19,1] = 'butorphanol tartrate 2 MG/ML Injectable Solution [Torbugesic]'
$data[19,2] = 'SBD'
$data[19,3] = 'Active'
$data[19,4] = '[''butorphanol'']'
$data[20,0] = 857192
$data[20,1] = 'butorphanol Injectable Solution [Dolorex Solution]'
$data[20,2] = 'SBDF'
$data[20,3] = 'Active'
$data[20,4] = '[''butorphanol'']'
$data[21,0] = 1307713
$data[21,1] = 'Kalopanax septemlobus bark extract'
$data[21,2] = 'IN'
$data[21,3] = 'Active'
$data[21,4] = '[''Kalopanax septemlobus bark extract'']'
$data[22,0] = 6468
$data[22,1] = 'loperamide'
$data[22,2] = 'IN'
$data[22,3] = 'Active'
$data[22,4] = '[''loperamide'']'
$data[23,0] = 1006892
$data[23,1] = 'belladonna alkaloids / kaolin / phenobarbital'
$data[23,2] = 'MIN'
$data[23,3] = 'Active'
$data[23,4] = '[''belladonna alkaloids'', ''kaolin'', ''phenobarbital'']'
$data[24,0] = 1007079
$data[24,1] = 'belladonna extract, USP / chlorpheniramine / phenylephrine / pyrilamine'
$data[24,2] = 'MIN'
$data[24,3] = 'Active'
$data[24,4] = '[''chlorpheniramine'', ''phenylephrine'', ''belladonna extract, USP'', ''pyrilamine'']'
$data[25,0] = 1007139
$data[25,1] = 'belladonna alkaloids / caffeine'
$data[25,2] = 'MIN'
$data[25,3] = 'Active'
$data[25,4] = '[''belladonna alkaloids'', ''caffeine'']'
$data[26,0] = 1007539
$data[26,1] = 'belladonna extract, USP / ephedrine'
$data[26,2] = 'MIN'
$data[26,3] = 'Active'
$data[26,4] = '[''ephedrine'', ''belladonna extract, USP'']'
$data[27,0] = 1007603
$data[27,1] = 'caffeine / ergotamine / levorotatory alkaloids of belladonna / pentobarbital'
$data[27,2] = 'MIN'
$data[27,3] = 'Active'
$data[27,4] = '[''caffeine'', ''levorotatory alkaloids of belladonna'', ''ergotamine'', ''pentobarbital'']'
$data[28,0] = 1007608
$data[28,1] = 'belladonna alkaloids / phenobarbital'
$data[28,2] = 'MIN'
$data[28,3] = 'Active'
$data[28,4] = '[''belladonna alkaloids'', ''phenobarbital'']'
$data[29,0] = 1007644
$data[29,1] = 'belladonna alkaloids / chlorpheniramine / phenylephrine / phenylpropanolamine'
$data[29,2] = 'MIN'
$data[29,3] = 'Active'
$data[29,4] = '[''belladonna alkaloids'', ''chlorpheniramine'', ''phenylephrine'', ''phenylpropanolamine'']'
$data[30,0] = 1007787
$data[30,1] = 'belladonna extract, USP / methenamine / salicylamide'
$data[30,2] = 'MIN'
$data[30,3] = 'Active'
$data[30,4] = '[''methenamine'', ''belladonna extract, USP'', ''salicylamide'']'
$data[31,0] = 1007893
$data[31,1] = 'belladonna alkaloids / caffeine / ergotamine / pentobarbital'
$data[31,2] = 'MIN'
$data[31,3] = 'Active'
$data[31,4] = '[''belladonna alkaloids'', ''caffeine'', ''ergotamine'', ''pentobarbital'']'
$data[32,0] = 1008045
$data[32,1] = 'belladonna extract, USP / phenobarbital'
$data[32,2] = 'MIN'
$data[32,3] = 'Active'
$data[32,4] = '[''phenobarbital'', ''belladonna extract, USP'']'
$data[33,0] = 1008287
$data[33,1] = 'belladonna extract, USP / charcoal'
$data[33,2] = 'MIN'
$data[33,3] = 'Active'
$data[33,4] = '[''charcoal'', ''belladonna extract, USP'']'
$data[34,0] = 1008529
$data[34,1] = 'belladonna extract, USP / chlorpheniramine / pheniramine / phenylpropanolamine'
$data[34,2] = 'MIN'
$data[34,3] = 'Active'
$data[34,4] = '[''chlorpheniramine'', ''pheniramine'', ''phenylpropanolamine'', ''belladonna extract, USP'']'
$data[35,0] = 1008547
$data[35,1] = 'belladonna extract, USP / butabarbital'
$data[35,2] = 'MIN'
$data[35,3] = 'Active'
$data[35,4] = '[''butabarbital'', ''belladonna extract, USP'']'
$data[36,0] = 1008619
$data[36,1] = 'belladonna alkaloids / caffeine / ergotamine / phenobarbital'
$data[36,2] = 'MIN'
$data[36,3] = 'Active'
$data[36,4] = '[''belladonna alkaloids'', ''caffeine'', ''ergotamine'', ''phenobarbital'']'
$data[37,0] = 1008938
$data[37,1] = 'belladonna alkaloids / simethicone'
$data[37,2] = 'MIN'
$data[37,3] = 'Active'
$data[37,4] = '[''belladonna alkaloids'', ''simethicone'']'
$data[38,0] = 1009102
$data[38,1] = 'ergotamine / levorotatory alkaloids of belladonna / phenobarbital'
$data[38,2] = 'MIN'
$data[38,3] = 'Active'
$data[38,4] = '[''levorotatory alkaloids of belladonna'', ''ergotamine'', ''phenobarbital'']'
$data[39,0] = 1043
$data[39,1] = 'apomorphine'
$data[39,2] = 'IN'
$data[39,3] = 'Active'
$data[39,4] = '[''apomorphine'']'
$data[40,0] = 1312380
$data[40,1] = 'Atropa belladonna fruiting top extract'
$data[40,2] = 'IN'
$data[40,3] = 'Active'
$data[40,4] = '[''Atropa belladonna fruiting top extract'']'
$data[41,0] = 1353220
$data[41,1] = 'Atropa belladonna root extract'
$data[41,2] = 'IN'
$data[41,3] = 'Active'
$data[41,4] = '[''Atropa belladonna root extract'']'
$data[42,0] = 1359
$data[42,1] = 'belladonna alkaloids'
$data[42,2] = 'IN'
$data[42,3] = 'Active'
$data[42,4] = '[''belladonna alkaloids'']'
$data[43,0] = 1363430
$data[43,1] = 'Atropa belladonna whole extract'
$data[43,2] = 'IN'
$data[43,3] = 'Active'
$data[43,4] = '[''Atropa belladonna whole extract'']'
$data[44,0] = 1811764
$data[44,1] = 'Atropa belladonna flowering top extract'
$data[44,2] = 'IN'
$data[44,3] = 'Active'
$data[44,4] = '[''Atropa belladonna flowering top extract'']'
$data[45,0] = 215451
$data[45,1] = 'aspirin / oxycodone hydrochloride / oxycodone terephthalate'
$data[45,2] = 'MIN'
$data[45,3] = 'Active'
$data[45,4] = '[''aspirin'']'
$data[46,0] = 221074
$data[46,1] = 'chlorpheniramine polistirex'
$data[46,2] = 'PIN'
$data[46,3] = 'Active'
$data[46,4] = '[''chlorpheniramine'']'
$data[47,0] = 221113
$data[47,1] = 'levorotatory alkaloids of belladonna'
$data[47,2] = 'IN'
$data[47,3] = 'Active'
$data[47,4] = '[''levorotatory alkaloids of belladonna'']'
$data[48,0] = 314517
$data[48,1] = 'belladonna leaf extract'
$data[48,2] = 'IN'
$data[48,3] = 'Active'
$data[48,4] = '[''belladonna leaf extract'']'
$data[49,0] = 544581
$data[49,1] = 'Apokyn'
$data[49,2] = 'BN'
$data[49,3] = 'Active'
$data[49,4] = '[''apomorphine'']'
$data[50,0] = 71225
$data[50,1] = 'apomorphine hydrochloride'
$data[50,2] = 'PIN'
$data[50,3] = 'Active'
$data[50,4] = '[''apomorphine'']'
$data[51,0] = 89781
$data[51,1] = 'belladonna extract, USP'
$data[51,2] = 'IN'
$data[51,3] = 'Active'
$data[51,4] = '[''belladonna extract, USP'']'
$data[52,0] = 42347
$data[52,1] = 'Bupropion'
$data[52,2] = 'IN'
$data[52,3] = 'Active'
$data[52,4] = '[''bupropion'']'

$ws.Range("A2:E54").Value = $data
